$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.879.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.37%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.901.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.59%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4594"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3809"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.73%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.60"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07722"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9811"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.91%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.927.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.961"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.672"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.63%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07073"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.80%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.006"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "84.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009533"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.004"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.846.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.337"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.97%  "

$ws.Range("E24").Value = "  -3.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.131.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.76%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.097"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.63%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.590"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.838"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.92%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09270"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8619"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.83%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.096"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.51%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.254"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.020"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.36%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05703"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.69%  "

$ws.Range("E38").Value = "  -2.53%  "

$ws.Range("E39").Value = "  +0.17%  "

$ws.Range("E40").Value = "  -4.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.472"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.70%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5507"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.80%  "

$ws.Range("E43").Value = "  -4.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.308"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.725"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5192"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.80%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.65%  "

$ws.Range("E48").Value = "  -4.70%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06832"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.93%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.69%  "

$ws.Range("E51").Value = "  -5.64%  "
